$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.192.48"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "3.063.90"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'561.93"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").Value = "'143.75"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.063.30"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").Value = "'6.13"
$ws.Range("E11").Value = "  -8.23%  "
$ws.Range("D12").Value = "'0.490"
$ws.Range("E12").Value = "  +10.67%  "
$ws.Range("E13").Value = "  +5.64%  "
$ws.Range("D14").Value = "'35.57"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "3.565.37"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "64.238.64"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").Value = "3.065.77"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("D20").Value = "'479.20"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").Value = "'13.95"
$ws.Range("E21").Value = "  +4.98%  "
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").Value = "'7.61"
$ws.Range("E23").Value = "  +5.90%  "
$ws.Range("D24").Value = "'14.42"
$ws.Range("E24").Value = "  +15.08%  "
$ws.Range("D25").Value = "'82.31"
$ws.Range("E25").Value = "  +4.19%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  +7.36%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "'26.38"
$ws.Range("E31").Value = "  +4.06%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("D34").Value = "'5.74"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("D35").Value = "'6.26"
$ws.Range("E35").Value = "  +7.85%  "
$ws.Range("D36").Value = "'54.97"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.0410"
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("D38").Value = "'447.27"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'0.0815"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").Value = "'2.85"
$ws.Range("E40").Value = "  +11.56%  "
$ws.Range("D41").Value = "3.012.74"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").Value = "'8.26"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'27.88"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("E45").Value = "  +7.17%  "
$ws.Range("E46").Value = "  +10.22%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").Value = "0.0₃0521"
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").Value = "'118.81"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  +4.69%  "
